# Added reward interpolation for mix of local/global reward
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new header column "O_reward_discount" right after the existing
# header in AA1, copying its style (bold header formatting).
$ws.Range("AA1").Copy() | Out-Null
$ws.Range("AB1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("AB1").Value = "O_reward_discount"

# Remove the sample data row (row 2) entirely.
$ws.Rows.Item(2).Delete()
